# ---------------------------------------------------------------------------
# Definitie.docx edit script
#
# 1) Replace the "Systeem:" intro paragraph: drop the _GoBack bookmark that
#    used to sit at its start, prepend two new sentences about the
#    application's OS requirements, and change "...gebruiken ze Microsoft
#    Excel om..." to "...gebruiken ROCit Microsoft Excel om..." (with ROCit
#    flagged as a proofing exception, matching the rest of the document).
# 2) Remove the now-superfluous leading space in "... hun verstuurde
#    kosten. ".
# 3) Re-add the _GoBack bookmark right after "€800.-".
# 4) Move <w:lastRenderedPageBreak/> from the "Advies:" run to the
#    "Tijdbesparing en geldbesparing." run.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Rewrite the "Systeem:" intro paragraph -----------------------------

$rng = $d.Content
$rng.Find.Execute("Op dit moment gebruiken ze Microsoft Excel om", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range

$paraInner = '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>De applicatie wordt gemaakt op</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> computers</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> die allemaal gebruik maken van </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>W</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>indows 10</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> als operating system</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>.</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> De applicatie kan op alle Operating systems </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">volledig worden gebruikt. </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Op dit moment ge</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">bruiken </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>ROCit</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Microsoft Excel om </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">de kosten van de onderhoudsbeurten </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>voor alle auto' + [char]0x2019 + 's bij te houden.</w:t></w:r>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p>' + $paraInner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$para.InsertXML($xml)

# --- 2) Fix the double space before "hun verstuurde kosten." ---------------

$d.Content.Find.Execute(" hun verstuurde kosten. ", $true, $false, $false, $false, $false, $true, 1, $false, `
  "hun verstuurde kosten. ", 2)

# --- 3) Re-create the _GoBack bookmark right after "(EUR)800.-" ------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$priceRng = $d.Content
$priceRng.Find.Execute([char]0x20AC + "800.-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$priceRng.Collapse(0)
# Insert a temporary marker character so the insertion point is no longer
# sitting exactly on the paragraph-end boundary (Bookmarks.Add mis-behaves
# there), anchor the bookmark around it, then delete the marker again - the
# bookmark collapses back onto the correct, now-safe position.
$priceRng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $priceRng)
$markerRng = $d.Range($priceRng.Start, $priceRng.End)
$markerRng.Text = ""

# --- 4) Move the lastRenderedPageBreak marker -------------------------------

# 4a) Drop it from the "Advies:" run (rewriting the paragraph without the
#     marker, but with identical visible text/formatting).
$advRng = $d.Content
$advRng.Find.Execute("Advies:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$advPara = $advRng.Paragraphs(1).Range

$advXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:lang w:val="nl-NL"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="nl-NL"/></w:rPr><w:t>Advies:</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$advPara.InsertXML($advXml)

# 4b) Add it to the "Tijdbesparing en geldbesparing." run.
$tijdRng = $d.Content
$tijdRng.Find.Execute("Tijdbesparing en geldbesparing.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tijdPara = $tijdRng.Paragraphs(1).Range

$tijdXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:lang w:val="nl-NL"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:lastRenderedPageBreak/><w:t>Tijdbesparing en geldbesparing.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tijdPara.InsertXML($tijdXml)

Write-Output "done"
